$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: person name + card number changes
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit-only string that must stay TEXT (not be
# auto-converted to a number) while keeping its original cell style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 30.05.2024"

# Row 6
$ws.Range("B6").Value = "01.06."
$ws.Range("C6").Value = "02.06."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 60674945"
$ws.Range("E6").Value = "84,01-"

# Row 7
$ws.Range("B7").Value = "02.06."
$ws.Range("C7").Value = "03.06."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 41007904"
$ws.Range("E7").Value = "42,40-"

# Row 8
$ws.Range("B8").Value = "06.06."
$ws.Range("C8").Value = "07.06."
$ws.Range("D8").Value = "PAYPAL BZKMDE"
$ws.Range("E8").Value = "17,25-"

# Row 9
$ws.Range("B9").Value = "10.06."
$ws.Range("C9").Value = "11.06."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-37655683"
$ws.Range("E9").Value = "57,60-"

# Rows 10 and 11 become empty (transactions removed)
$ws.Range("B10:E11").Value = ""
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 15.06.2024"
$ws.Range("E12").Value = "201,26-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 24.06.2024"
